$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.873.52"

$ws.Range("D3").Value = "1.643.86"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("E6").Value = "  -0.47%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.99"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.45%  "

$ws.Range("E9").Value = "  +0.84%  "

$ws.Range("E10").Value = "  -0.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0900"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.61%  "

$ws.Range("D12").Value = "1.874.91"
$ws.Range("E12").Value = "  +0.91%  "

$ws.Range("D13").Value = "1.639.70"
$ws.Range("E13").Value = "  +0.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.595"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.71"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +9.14%  "

$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("D17").Value = "29.861.15"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.45"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "238.28"
$ws.Range("D19").ClearFormats()

$ws.Range("D20").Value = "0.0₃0705"
$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.95"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.14"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.17"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.92"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.65"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.68"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0497"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.50%  "

$ws.Range("E31").Value = "  -0.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.39"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.20"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.03%  "

$ws.Range("D34").Value = "1.422.11"
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("E35").Value = "  +2.88%  "

$ws.Range("E36").Value = "  -0.79%  "

$ws.Range("E37").Value = "  +2.04%  "

$ws.Range("E38").Value = "  -7.20%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  +2.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "76.69"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +10.54%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.838"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.60%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0503"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.94"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "50.63"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.33%  "

$ws.Range("D48").Value = "1.783.10"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("E49").Value = "  -1.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "93.95"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.74%  "

$ws.Range("D51").Value = "0.0₆0110"
$ws.Range("E51").Value = "  +1.88%  "
